$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metodologia")

# --- Row 3-10 block (Table35 top section, columns I/J/K = Springler/Wiley/Taylor & Francis) ---
$ws.Range("J3").Value = 5440
$ws.Range("K3").Value = 4610

$ws.Range("J4").Value = 27414
$ws.Range("K4").Value = 11745

$ws.Range("J5").Value = 69147
$ws.Range("K5").Value = 32024

$ws.Range("J6").Value = 4961
$ws.Range("K6").Value = 4306

$ws.Range("I7").Value = 3453
$ws.Range("J7").Value = 773
$ws.Range("K7").Value = 1559

$ws.Range("I8").Value = 2506
$ws.Range("J8").Value = 238

$ws.Range("I9").Value = 94925
$ws.Range("J9").Value = 13442

$ws.Range("I10").Value = 308
$ws.Range("J10").Value = 59

# --- Row 14-21 block (Table35 bottom / Table33, columns I/J/K and Q/R/S) ---
$ws.Range("J14").Value = 2950
$ws.Range("K14").Value = 3289
$ws.Range("R14").Value = 40
$ws.Range("S14").Value = 70

$ws.Range("J15").Value = 12553
$ws.Range("K15").Value = 7383
$ws.Range("R15").Value = 60
$ws.Range("S15").Value = 50

$ws.Range("J16").Value = 31700
$ws.Range("K16").Value = 19240
$ws.Range("R16").Value = 20
$ws.Range("S16").Value = 20

$ws.Range("J17").Value = 2674
$ws.Range("K17").Value = 3055
$ws.Range("Q17").Value = 80
$ws.Range("R17").Value = 60
$ws.Range("S17").Value = 60

$ws.Range("I18").Value = 1459
$ws.Range("J18").Value = 411
$ws.Range("K18").Value = 1145
$ws.Range("Q18").Value = 100
$ws.Range("R18").Value = 20
$ws.Range("S18").Value = 70

$ws.Range("I19").Value = 1041
$ws.Range("J19").Value = 114
$ws.Range("Q19").Value = 80
$ws.Range("R19").Value = 20

$ws.Range("I20").Value = 28492
$ws.Range("J20").Value = 5747
$ws.Range("Q20").Value = 20
$ws.Range("R20").Value = 60

$ws.Range("I21").Value = 154
$ws.Range("J21").Value = 40
$ws.Range("Q21").Value = 40
$ws.Range("R21").Value = 40

# --- Update the view to match where the author ended up working ---
$ws.Activate()
$ws.Range("M18").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 12
